$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '74.857.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.67%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.819.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +7.09%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '187.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.99%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '594.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.01%  '

# Row 8
$ws.Range('E8').Value = '  +2.87%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.192'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.17%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.816.75'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.09%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.161'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.30%  '

# Row 12
$ws.Range('E12').Value = '  +3.55%  '

# Row 13
$ws.Range('E13').Value = '  +2.75%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.336.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.19%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '74.781.62'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.82%  '

# Row 16
$ws.Range('E16').Value = '  -1.49%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.01%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.819.35'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.13%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.93'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.04%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.74%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '377.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.38%  '

# Row 22
$ws.Range('E22').Value = '  -1.70%  '

# Row 23
$ws.Range('E23').Value = '  -0.83%  '

# Row 24
$ws.Range('E24').Value = '  -0.13%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.86'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.05%  '

# Row 26
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.81'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.68%  '

# Row 27
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.965.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.23%  '

# Row 28
$ws.Range('B28').Value = 'NEARProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.13%  '

# Row 29
$ws.Range('E29').Value = '  +9.61%  '

# Row 30
$ws.Range('E30').Value = '  -0.10%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '517.09'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.02%  '

# Row 32
$ws.Range('E32').Value = '  -0.45%  '

# Row 33
$ws.Range('E33').Value = '  +0.45%  '

# Row 34
$ws.Range('E34').Value = '  +2.51%  '

# Row 35
$ws.Range('E35').Value = '  -0.02%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '163.56'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.24%  '

# Row 37
$ws.Range('E37').Value = '  +3.88%  '

# Row 38
$ws.Range('E38').Value = '  -1.28%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.37'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.43%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '186.56'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +16.03%  '

# Row 41
$ws.Range('E41').Value = '  -0.01%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.340'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.53%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.53%  '

# Row 44
$ws.Range('E44').Value = '  -0.44%  '

# Row 45
$ws.Range('E45').Value = '  +1.67%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.76%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0857'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.28%  '

# Row 48
$ws.Range('E48').Value = '  -2.47%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.578'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.16%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.70'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.26%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.635'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.10%  '
